# Change gapfilling indexing system
# Berge_MDS sheet: shift the Proxy_vars / Proxy_vars_subset / Proxy_vars_range
# values so a new proxy variable ("delta_Tair_Teau") becomes the first entry,
# the old entries shift down one row, and the last (4th) numeric range value
# wraps around to become the new first one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berge_MDS")

# --- Row 4: air_temp_HMP45C -> air_relativeHumidity, range 3.5 -> 3 -------
$ws.Range("D4").Value = "air_relativeHumidity"
$ws.Range("F4").Value = 3

# --- Row 3: air_relativeHumidity -> wind_speed_05103, E3 cleared, range 3->1
$ws.Range("D3").Value = "wind_speed_05103"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 1

# --- Row 2: new proxy var "delta_Tair_Teau", range 1 -> 3.5 ---------------
$ws.Range("D2").Value = "delta_Tair_Teau"
$ws.Range("E2").Value = "delta_Tair_Teau"
$ws.Range("F2").Value = 3.5

# D2/E2 lose their explicit (no-op) cell style in the target file, while D4
# gains the style that D3 already carries.
$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearFormats()

# Selection cursor moved from E17 to F17 in the saved view state.
$ws.Range("F17").Select()
